$d = $word.ActiveDocument

# --- Body content -----------------------------------------------------
# Paragraph 2 ("юдьододщожщжщжз") is removed entirely.
$d.Paragraphs(2).Range.Delete()

# Paragraph 1 ("4515446465") becomes a single "k" at 72pt (sz 144 half-points).
$p1 = $d.Paragraphs(1)
$p1.Range.Text = "k"
$r1 = $d.Range($p1.Range.Start, $p1.Range.Start + 1)
$r1.Font.Size = 72

# Remaining new paragraphs to append, in order. The last one also gets an
# explicit black font color.
$texts = @("k", "k", "kk", "k", "k", "k", "k", "kkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkkk")

$lastRange = $d.Paragraphs(1).Range
foreach ($t in $texts) {
    $lastRange.InsertParagraphAfter() | Out-Null
    $lastRange = $d.Paragraphs($d.Paragraphs.Count).Range
    $lastRange.Text = $t
    $charStart = $lastRange.Start
    $charEnd = $charStart + $t.Length
    $fr = $d.Range($charStart, $charEnd)
    $fr.Font.Size = 72
}

# Final paragraph ("kkkk...") additionally carries an explicit black color.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastFullRange = $d.Range($lastPara.Range.Start, $lastPara.Range.Start + $texts[$texts.Length - 1].Length)
$lastFullRange.Font.Color = 0

# --- Styles -------------------------------------------------------------
# Normal style now carries an explicit Calibri rFonts declaration.
$normal = $d.Styles("Normal")
$normal.Font.Name = "Calibri"
$normal.Font.NameFarEast = "Calibri"
